# This script re-applies a batch of stock-ledger corrections to the
# "CryCompanywiseStockReport" sheet. For a number of duplicate stock-item
# rows (same item code/description, two receipt batches with different
# rates), the Code/Rate2/Qty/Value columns (B, D, E, F, G) had been
# entered against the wrong batch row and need to be swapped back onto
# the correct row. A couple of single (non-duplicated) rows also get a
# straight quantity/value correction, and the Sub Total / Grand Total
# rows are refreshed to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Row($RowA, $RowB, $Columns) {
    foreach ($col in $Columns) {
        $refA = "$col$RowA"
        $refB = "$col$RowB"
        $valA = $ws.Range($refA).Value2
        $valB = $ws.Range($refB).Value2
        $ws.Range($refA).Value2 = $valB
        $ws.Range($refB).Value2 = $valA
    }
}

$cols = @("B", "D", "E", "F", "G")

# Pairs of duplicate item rows whose batch data (Code/Rate2/Qty/Value)
# were swapped back to the correct rows.
$pairs = @(
    @(136, 137),
    @(246, 247),
    @(292, 293),
    @(311, 312),
    @(420, 421),
    @(472, 473),
    @(476, 477),
    @(479, 480),
    @(485, 486),
    @(564, 565),
    @(596, 597),
    @(705, 706),
    @(732, 733)
)

foreach ($p in $pairs) {
    Swap-Row $p[0] $p[1] $cols
}

# Rows 294/295/296 form a 3-way rotation rather than a simple swap:
# new294 <- old295, new295 <- old296, new296 <- old294
foreach ($col in $cols) {
    $ref294 = "${col}294"
    $ref295 = "${col}295"
    $ref296 = "${col}296"
    $v294 = $ws.Range($ref294).Value2
    $v295 = $ws.Range($ref295).Value2
    $v296 = $ws.Range($ref296).Value2
    $ws.Range($ref294).Value2 = $v295
    $ws.Range($ref295).Value2 = $v296
    $ws.Range($ref296).Value2 = $v294
}

# Single (non-duplicated) item rows with a direct quantity/value correction.
$ws.Range("F313").Value2 = 23
$ws.Range("G313").Value2 = 2791.28

$ws.Range("F482").Value2 = 437
$ws.Range("G482").Value2 = 2875.46

# Refresh the Sub Total / Grand Total rows affected by the above corrections.
$ws.Range("B339").Value2 = 345188.55
$ws.Range("B492").Value2 = 1476.55
$ws.Range("B793").Value2 = 3238411.15
$ws.Range("B794").Value2 = 3238411.15
